$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new dataset entry "Flowise" to the list in column A
$ws.Range("A14").Value = "Flowise"

# Move the active selection to mirror the saved view state
$ws.Range("G18").Select()
